$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels and Devices")

# --- Copy row 8's formatting down into the new row 9 first, so every
#     cell in row 9 picks up the same styles used by row 8 (s=6,6,10,6,6,13,6,6,6,3,3,14)
$ws.Range("A8:L8").Copy()
$ws.Range("A9:L9").PasteSpecial(-4122)

# --- Update the existing row 8 values that changed
#     F8 / K8 "Expected 24V PSU Load value" / "Expected 3rd 24V PSU Load value" text 0.329 -> 0.319
#     J8 "Expected 2nd 24V PSU Load value" 0.405 -> 0.395
$ws.Range("F8").Value = "'0.319"
$ws.Range("J8").Value = 0.395
$ws.Range("K8").Value = 0.319

# --- Populate the newly added row 9 (new loop-card device entry)
$ws.Range("A9").Value = "MZX252"
$ws.Range("B9").Value = "Node1"
$ws.Range("D9").Value = "PFI"
$ws.Range("E9").Value = 14
$ws.Range("F9").Value = "'0.100"
$ws.Range("H9").Value = "XLM800-Zetfas-C"
$ws.Range("G9").Value = "XLM800-Zetfas"
$ws.Range("I9").Value = "Loops"
$ws.Range("J9").Value = 0.165
$ws.Range("K9").Value = "'0.100"
$ws.Range("L9").Value = "24V Rail(A)"

# --- Move the active selection to B9 (was L8), and drop any frozen/scrolled
#     top-left-cell state so the sheet view matches a fresh selection on B9.
$ws.Range("B9").Select() | Out-Null
